# "los tipos de datos estan afectando todo pa, terrible." --
# add the missing rows to Clientes and Productos, and fix the
# "Marca" column width on Productos.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Clientes: 3 new customer rows (A:C only -- "Activo" left blank)
# ---------------------------------------------------------------
$clientes = $wb.Worksheets.Item("Clientes")

$clientesRows = @(
    @(1000535410, "Albeiro Molina", 3245619849),
    @(1888473247, "Pedro la piedra", 2342345467),
    @(1222323423, "Validasras", 323234323)
)

$r = 3
foreach ($row in $clientesRows) {
    $clientes.Cells.Item($r, 1).Value = $row[0]
    $clientes.Cells.Item($r, 2).Value = $row[1]
    $clientes.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Productos: 5 new product rows (A:H) + widen the "Marca" column
# ---------------------------------------------------------------
$productos = $wb.Worksheets.Item("Productos")

$productosRows = @(
    @("asdfsdfa",  2423423423423, "sfsfsdf", 324234,  324234,  234234,  $true, "04/06/2024 02:49"),
    @("edadasdad", 2313131232131, "fafasf",  223423,  223434,  232423,  $true, "04/06/2024 03:00"),
    @("dsfsdf",    2313131232133, "fafasf",  223423,  223434,  232423,  $true, "04/06/2024 03:02"),
    @("teclados",  3245234543252, "Basura",  34252,   324534,  245234,  $true, "04/06/2024 04:44"),
    @("sfdsfafa",  2342342343242, "wqqtrre", 23424,   3424234, 2342342, $true, "04/06/2024 04:46")
)

$r = 9
foreach ($row in $productosRows) {
    $productos.Cells.Item($r, 1).Value = $row[0]
    $productos.Cells.Item($r, 2).Value = $row[1]
    $productos.Cells.Item($r, 3).Value = $row[2]
    $productos.Cells.Item($r, 4).Value = $row[3]
    $productos.Cells.Item($r, 5).Value = $row[4]
    $productos.Cells.Item($r, 6).Value = $row[5]
    $productos.Cells.Item($r, 7).Value = $row[6]
    $productos.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Column C ("Marca") grows from width 8 to width 9. Excel pads the
# character width it stores in the XML by 5/6 above whatever is
# assigned to ColumnWidth, so back that padding out here.
$productos.Columns.Item(3).ColumnWidth = 9 - (5 / 6)
